$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.113.55'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '1.957.20'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.20'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4909'
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2987'
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06857'
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.24'
$ws.Range("E10").Value = '  -1.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '108.41'
$ws.Range("E11").Value = '  -3.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07768'
$ws.Range("E12").Value = '  +1.61%  '
$ws.Range("D13").Value = '1.930.70'
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.466'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7111'
$ws.Range("E15").Value = '  +2.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.33'
$ws.Range("E16").Value = '  -4.07%  '
$ws.Range("D17").Value = '31.141.55'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.26'
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007764'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").Value = '2.179.30'
$ws.Range("E21").Value = '  -1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.498'
$ws.Range("E22").Value = '  -3.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.537'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.816'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.64'
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.11'
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.226'
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1054'
$ws.Range("E29").Value = '  -3.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.433'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.584'
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.593'
$ws.Range("E32").Value = '  -3.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.452'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04975'
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7591'
$ws.Range("E35").Value = '  -2.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.182'
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.737'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02041'
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.709'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.178'
$ws.Range("E40").Value = '  +6.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.451'
$ws.Range("E41").Value = '  +8.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4513'
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.279'
$ws.Range("E43").Value = '  +12.24%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '73.38'
$ws.Range("E44").Value = '  +3.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.51'
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.8834'
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.459'
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '962.82'
$ws.Range("E49").Value = '  +6.14%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1267'
$ws.Range("E50").Value = '  +1.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2598'
$ws.Range("E51").Value = '  +2.06%  '
